# repull data, push all data, mean calculation
# Updates the dSF column (F) values for each game row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -6
    3  = -2
    4  = 6
    5  = 5
    7  = -1
    8  = -2
    9  = 5
    11 = 4
    12 = 4
    13 = 2
    14 = -2
    16 = 1
    18 = 4
    20 = -1
    21 = 2
    23 = 2
    24 = 3
    25 = 3
    26 = 4
    27 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
